$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: record A=111863073
$ws.Range("A2").Value = 111863073
$ws.Range("B2").Value = 89033
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 3286
$ws.Range("F2").Value = "Flattoppad klubbsvamp"
$ws.Range("G2").Value = "Clavariadelphus truncatus"
$ws.Range("H2").Value = "(Quél.) Donk"
$ws.Range("I2").Value = "2"
$ws.Range("J2").Value = "fruktkroppar"
$ws.Range("P2").Value = "Charlottenberg (Charlottenberg), Upl"
$ws.Range("Q2").Value = 655228
$ws.Range("R2").Value = 6634879
$ws.Range("Z2").Value = "10:50"
$ws.Range("AB2").Value = "10:50"
$ws.Range("AC2").ClearContents()

# Row 3: record A=111863045
$ws.Range("A3").Value = 111863045
$ws.Range("B3").Value = 89033
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 3286
$ws.Range("F3").Value = "Flattoppad klubbsvamp"
$ws.Range("G3").Value = "Clavariadelphus truncatus"
$ws.Range("H3").Value = "(Quél.) Donk"
$ws.Range("I3").Value = "11"
$ws.Range("J3").Value = "fruktkroppar"
$ws.Range("P3").Value = "Charlottenberg (Charlottenberg), Upl"
$ws.Range("Q3").Value = 655234
$ws.Range("R3").Value = 6634889
$ws.Range("Z3").Value = "10:50"
$ws.Range("AB3").Value = "10:50"
$ws.Range("AC3").Value = "Under gran i svacka"

# Row 4: record A=111862959
$ws.Range("A4").Value = 111862959
$ws.Range("B4").Value = 90821
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 5964
$ws.Range("F4").Value = "Fjällig taggsvamp s.str."
$ws.Range("G4").Value = "Sarcodon imbricatus s.str."
$ws.Range("H4").Value = "(L.:Fr.) P.Karst."
$ws.Range("I4").Value = "11"
$ws.Range("J4").Value = "fruktkroppar"
$ws.Range("P4").Value = "Charlottenberg, Upl"
$ws.Range("Q4").Value = 655218
$ws.Range("R4").Value = 6634940
$ws.Range("Z4").Value = "10:37"
$ws.Range("AB4").Value = "10:37"
$ws.Range("AC4").Value = "Under gran och tall i en svacka"

# Row 5: record A=111863001
$ws.Range("A5").Value = 111863001
$ws.Range("B5").Value = 90466
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 4769
$ws.Range("F5").Value = "Svavelriska"
$ws.Range("G5").Value = "Lactarius scrobiculatus"
$ws.Range("H5").Value = "(Scop.:Fr.) Fr."
$ws.Range("I5").Value = "1"
$ws.Range("J5").Value = "fruktkroppar"
$ws.Range("P5").Value = "Charlottenberg, Upl"
$ws.Range("Q5").Value = 655218
$ws.Range("R5").Value = 6634940
$ws.Range("Z5").Value = "10:47"
$ws.Range("AB5").Value = "10:47"
$ws.Range("AC5").ClearContents()

# Row 6: record A=111863269
$ws.Range("A6").Value = 111863269
$ws.Range("B6").Value = 85183
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 249278
$ws.Range("F6").Value = "Barrviolspindling"
$ws.Range("G6").Value = "Cortinarius harcynicus"
$ws.Range("H6").Value = "(Pers.) M.M.Moser"
$ws.Range("I6").Value = "4"
$ws.Range("J6").Value = "fruktkroppar"
$ws.Range("P6").Value = "Charlottenberg (Charlottenberg), Upl"
$ws.Range("Q6").Value = 655135
$ws.Range("R6").Value = 6634800
$ws.Range("Z6").Value = "11:02"
$ws.Range("AB6").Value = "11:02"
$ws.Range("AC6").Value = "4 ex i gräsglänta under gran och tall."

# Row 7: record A=111863218
$ws.Range("A7").Value = 111863218
$ws.Range("B7").Value = 90155
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 6031
$ws.Range("F7").Value = "Blomkålssvamp"
$ws.Range("G7").Value = "Sparassis crispa"
$ws.Range("H7").Value = "(Wulfen:Fr.) Fr."
$ws.Range("I7").Value = "1"
$ws.Range("J7").Value = "fruktkroppar"
$ws.Range("P7").Value = "Charlottenberg (Charlottenberg), Upl"
$ws.Range("Q7").Value = 655138
$ws.Range("R7").Value = 6634821
$ws.Range("Z7").Value = "10:53"
$ws.Range("AB7").Value = "10:53"
$ws.Range("AC7").ClearContents()

# Row 8: record A=111863288
$ws.Range("A8").Value = 111863288
$ws.Range("B8").Value = 85183
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 249278
$ws.Range("F8").Value = "Barrviolspindling"
$ws.Range("G8").Value = "Cortinarius harcynicus"
$ws.Range("H8").Value = "(Pers.) M.M.Moser"
$ws.Range("I8").Value = "1"
$ws.Range("J8").Value = "fruktkroppar"
$ws.Range("P8").Value = "Charlottenberg (Charlottenberg), Upl"
$ws.Range("Q8").Value = 655135
$ws.Range("R8").Value = 6634793
$ws.Range("Z8").Value = "11:02"
$ws.Range("AB8").Value = "11:02"
$ws.Range("AC8").ClearContents()

# Row 9: record A=111863040
$ws.Range("A9").Value = 111863040
$ws.Range("B9").Value = 90821
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 5964
$ws.Range("F9").Value = "Fjällig taggsvamp s.str."
$ws.Range("G9").Value = "Sarcodon imbricatus s.str."
$ws.Range("H9").Value = "(L.:Fr.) P.Karst."
$ws.Range("I9").Value = ""
$ws.Range("J9").ClearContents()
$ws.Range("P9").Value = "Charlottenberg (Charlottenberg), Upl"
$ws.Range("Q9").Value = 655235
$ws.Range("R9").Value = 6634878
$ws.Range("Z9").Value = "10:49"
$ws.Range("AB9").Value = "10:49"
$ws.Range("AC9").Value = "Halv häxring, 3 m i diameter"

# Row 10: record A=111863402
$ws.Range("A10").Value = 111863402
$ws.Range("B10").Value = 90821
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 5964
$ws.Range("F10").Value = "Fjällig taggsvamp s.str."
$ws.Range("G10").Value = "Sarcodon imbricatus s.str."
$ws.Range("H10").Value = "(L.:Fr.) P.Karst."
$ws.Range("I10").Value = "1"
$ws.Range("J10").Value = "fruktkroppar"
$ws.Range("P10").Value = "Charlottenberg (Charlottenberg), Upl"
$ws.Range("Q10").Value = 655200
$ws.Range("R10").Value = 6634770
$ws.Range("Z10").Value = "11:02"
$ws.Range("AB10").Value = "11:02"
$ws.Range("AC10").ClearContents()
